# Updated final report formatting
#
# 1. Bump the cached "datetimeFigureOut" field text (4/18/2024 -> 4/19/2024)
#    on the slide master and every slide layout's Date Placeholder.
# 2. Reposition/resize a handful of callout text boxes on slides 5 and 8.
# 3. Reposition/resize and restyle (solid fill -> outline only) the two
#    right-arrow callouts on slide 9.

function EmuToPt($emu) {
    # The host floors fractional points when converting back to EMU, so
    # nudge by 0.5 EMU worth of points to land exactly on the target EMU.
    return ([double]$emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text: slide master + all custom layouts
# ---------------------------------------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "4/18/2024") {
            $sh.TextFrame.TextRange.Text = "4/19/2024"
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "4/18/2024") {
                $sh.TextFrame.TextRange.Text = "4/19/2024"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 5 ("Database Schema") - TextBox 39 callout
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$tb39 = $slide5.Shapes.Item("TextBox 39")
$tb39.Left = EmuToPt 8671352
$tb39.Top = EmuToPt 4636499
$tb39.Width = EmuToPt 2552076
$tb39.Height = EmuToPt 1021556

# ---------------------------------------------------------------------
# 3. Slide 8 ("Modeling Flow") - three callout text boxes
# ---------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)

$tb70 = $slide8.Shapes.Item("TextBox 70")
$tb70.Left = EmuToPt 831301
$tb70.Top = EmuToPt 4600641
$tb70.Width = EmuToPt 2623459
$tb70.Height = EmuToPt 923330

$tb71 = $slide8.Shapes.Item("TextBox 71")
$tb71.Left = EmuToPt 828806
$tb71.Top = EmuToPt 5512822
$tb71.Width = EmuToPt 2623459
$tb71.Height = EmuToPt 923330

$tb73 = $slide8.Shapes.Item("TextBox 73")
$tb73.Left = EmuToPt 905601
$tb73.Top = EmuToPt 3237540
$tb73.Width = EmuToPt 2437674
$tb73.Height = EmuToPt 1021556

# ---------------------------------------------------------------------
# 4. Slide 9 ("Testing & Evaluation") - the two right-arrow callouts
#    move/shrink slightly and swap from a solid fill to an outline only.
# ---------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)

$arrow22 = $slide9.Shapes.Item("Right Arrow 22")
$arrow22.Left = EmuToPt 3924299
$arrow22.Top = EmuToPt 3166332
$arrow22.Width = EmuToPt 767769
$arrow22.Height = EmuToPt 628377
$arrow22.Fill.Visible = 0
$arrow22.Line.Visible = -1
$arrow22.Line.ForeColor.RGB = 15778605

$arrow24 = $slide9.Shapes.Item("Right Arrow 24")
$arrow24.Left = EmuToPt 6987612
$arrow24.Top = EmuToPt 3163471
$arrow24.Width = EmuToPt 767769
$arrow24.Height = EmuToPt 628377
$arrow24.Fill.Visible = 0
$arrow24.Line.Visible = -1
$arrow24.Line.ForeColor.RGB = 15778605
